# Loan RBI, Variable Instalments
#
# On the "Repayment Schedule" sheet, insert a new (blank) column before
# column N. This shifts the existing "Late" column (N -> O) and the
# existing "Outstanding" column (P -> Q) one place to the right, leaving
# a new empty column N in between them.
#
# Also make "Repayment Schedule" the active sheet/tab (it was
# "Transactions" before) and update the active selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Switch focus to the Repayment Schedule tab (was Transactions).
$ws.Activate()

# Insert a blank column before N - shifts N (Late) -> O, P (Outstanding) -> Q.
$ws.Columns("N").Insert()

# Update the active cell/selection on the Repayment Schedule sheet.
$ws.Range("T6").Select()
